{"js": "// Replace the date line and each of the 25 multiplication-table entries\n// with their updated values, per the commit's regenerated worksheet data.\nconst replacements = [\n  [\"2024-05-17 Friday\", \"2024-05-18 Saturday\"],\n  [\"698\u00d75=\", \"925\u00d72=\"],\n  [\"521\u00d72=\", \"890\u00d72=\"],\n  [\"337\u00d78=\", \"371\u00d72=\"],\n  [\"779\u00d77=\", \"978\u00d77=\"],\n  [\"319\u00d76=\", \"471\u00d76=\"],\n  [\"855\u00d79=\", \"461\u00d73=\"],\n  [\"506\u00d76=\", \"112\u00d78=\"],\n  [\"134\u00d77=\", \"837\u00d72=\"],\n  [\"120\u00d75=\", \"410\u00d78=\"],\n  [\"700\u00d76=\", \"453\u00d78=\"],\n  [\"743\u00d76=\", \"154\u00d72=\"],\n  [\"734\u00d78=\", \"480\u00d78=\"],\n  [\"344\u00d73=\", \"660\u00d78=\"],\n  [\"634\u00d77=\", \"566\u00d74=\"],\n  [\"663\u00d74=\", \"639\u00d76=\"],\n  [\"148\u00d79=\", \"330\u00d72=\"],\n  [\"560\u00d78=\", \"233\u00d79=\"],\n  [\"134\u00d78=\", \"178\u00d76=\"],\n  [\"327\u00d74=\", \"701\u00d79=\"],\n  [\"542\u00d75=\", \"190\u00d76=\"],\n  [\"480\u00d74=\", \"262\u00d76=\"],\n  [\"292\u00d79=\", \"767\u00d76=\"],\n  [\"295\u00d72=\", \"819\u00d75=\"],\n  [\"807\u00d79=\", \"910\u00d79=\"],\n  [\"876\u00d74=\", \"347\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the 25 multiplication-table entries\n# with their updated values, per the commit's regenerated worksheet data.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-17 Friday\", \"2024-05-18 Saturday\"),\n    @(\"698\u00d75=\", \"925\u00d72=\"),\n    @(\"521\u00d72=\", \"890\u00d72=\"),\n    @(\"337\u00d78=\", \"371\u00d72=\"),\n    @(\"779\u00d77=\", \"978\u00d77=\"),\n    @(\"319\u00d76=\", \"471\u00d76=\"),\n    @(\"855\u00d79=\", \"461\u00d73=\"),\n    @(\"506\u00d76=\", \"112\u00d78=\"),\n    @(\"134\u00d77=\", \"837\u00d72=\"),\n    @(\"120\u00d75=\", \"410\u00d78=\"),\n    @(\"700\u00d76=\", \"453\u00d78=\"),\n    @(\"743\u00d76=\", \"154\u00d72=\"),\n    @(\"734\u00d78=\", \"480\u00d78=\"),\n    @(\"344\u00d73=\", \"660\u00d78=\"),\n    @(\"634\u00d77=\", \"566\u00d74=\"),\n    @(\"663\u00d74=\", \"639\u00d76=\"),\n    @(\"148\u00d79=\", \"330\u00d72=\"),\n    @(\"560\u00d78=\", \"233\u00d79=\"),\n    @(\"134\u00d78=\", \"178\u00d76=\"),\n    @(\"327\u00d74=\", \"701\u00d79=\"),\n    @(\"542\u00d75=\", \"190\u00d76=\"),\n    @(\"480\u00d74=\", \"262\u00d76=\"),\n    @(\"292\u00d79=\", \"767\u00d76=\"),\n    @(\"295\u00d72=\", \"819\u00d75=\"),\n    @(\"807\u00d79=\", \"910\u00d79=\"),\n    @(\"876\u00d74=\", \"347\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
